# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G (header "K") values for rows 2-38 change; everything else
# on the sheet stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value for column G
$newK = @{
    2  = 3
    3  = 0
    4  = 9
    5  = 6
    6  = 3
    7  = 7
    8  = 3
    9  = 3
    10 = 3
    11 = 4
    12 = 6
    13 = 2
    14 = 6
    15 = 4
    16 = 4
    17 = 4
    18 = 5
    19 = 4
    20 = 1
    21 = 9
    22 = 3
    23 = 6
    24 = 9
    25 = 4
    26 = 6
    27 = 5
    28 = 5
    29 = 3
    30 = 3
    31 = 7
    32 = 5
    33 = 5
    34 = 7
    35 = 1
    36 = 1
    37 = 2
    38 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
